# Update 'want to go' counts (column F) across sheets per the source data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4928  # was 4916
$ws.Range("F3").Value = 2770  # was 2768
$ws.Range("F5").Value = 2831  # was 2829
$ws.Range("F9").Value = 1742  # was 1741
$ws.Range("F11").Value = 492  # was 490
$ws.Range("F12").Value = 247  # was 245
$ws.Range("F13").Value = 430  # was 427
$ws.Range("F14").Value = 1078  # was 1077
$ws.Range("F15").Value = 303  # was 302
$ws.Range("F18").Value = 79  # was 78
$ws.Range("F19").Value = 1044  # was 1023
$ws.Range("F22").Value = 669  # was 667
$ws.Range("F23").Value = 754  # was 752
$ws.Range("F24").Value = 152  # was 150
$ws.Range("F25").Value = 11  # was 10
$ws.Range("F27").Value = 547  # was 538
$ws.Range("F28").Value = 48  # was 45
$ws.Range("F29").Value = 1662  # was 1661
$ws.Range("F30").Value = 1638  # was 1616
$ws.Range("F31").Value = 394  # was 385
$ws.Range("F32").Value = 46  # was 45
$ws.Range("F33").Value = 1557  # was 1545
$ws.Range("F34").Value = 219  # was 214
$ws.Range("F35").Value = 2386  # was 2378
$ws.Range("F36").Value = 408  # was 404
$ws.Range("F38").Value = 625  # was 623
$ws.Range("F40").Value = 68  # was 67
$ws.Range("F42").Value = 813  # was 807
$ws.Range("F43").Value = 1507  # was 1503
$ws.Range("F44").Value = 227  # was 223
$ws.Range("F46").Value = 502  # was 499
$ws.Range("F47").Value = 59  # was 57
$ws.Range("F49").Value = 117  # was 114

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 104  # was 101
$ws.Range("F12").Value = 46  # was 44

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4928  # was 4916
$ws.Range("F3").Value = 2770  # was 2768
$ws.Range("F4").Value = 2831  # was 2829
$ws.Range("F5").Value = 1742  # was 1741
$ws.Range("F9").Value = 492  # was 490
$ws.Range("F10").Value = 247  # was 245
$ws.Range("F11").Value = 430  # was 427
$ws.Range("F12").Value = 1078  # was 1077
$ws.Range("F13").Value = 303  # was 302
$ws.Range("F15").Value = 1044  # was 1024
$ws.Range("F17").Value = 669  # was 667
$ws.Range("F18").Value = 754  # was 752
$ws.Range("F19").Value = 152  # was 150
$ws.Range("F20").Value = 104  # was 101
$ws.Range("F21").Value = 104  # was 101
$ws.Range("F23").Value = 11  # was 10
$ws.Range("F26").Value = 547  # was 538
$ws.Range("F27").Value = 1662  # was 1661
$ws.Range("F28").Value = 1638  # was 1616
$ws.Range("F29").Value = 394  # was 385
$ws.Range("F30").Value = 46  # was 45
$ws.Range("F33").Value = 2386  # was 2378
$ws.Range("F34").Value = 408  # was 404
$ws.Range("F39").Value = 46  # was 44
$ws.Range("F41").Value = 68  # was 67
$ws.Range("F43").Value = 813  # was 807
$ws.Range("F44").Value = 1507  # was 1503
$ws.Range("F46").Value = 227  # was 223
$ws.Range("F47").Value = 502  # was 499
$ws.Range("F48").Value = 59  # was 57
